# Remove the "biosat" and "O2_Ar_ratio" attribute rows from the NCP output
# metadata sheet (ColumnHeaders). These were rows 7 and 8; deleting them
# shifts "ncp" and "k" up to rows 7 and 8 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnHeaders")

# Delete row 7 (biosat) and what was row 8 (O2_Ar_ratio, now also row 7
# after the first delete) as a single contiguous range so both rows are
# removed and subsequent rows (ncp, k) shift upward.
$ws.Range("A7:G8").EntireRow.Delete() | Out-Null

# Update the active selection to match the post-edit state.
$ws.Range("B16").Select() | Out-Null
